# "ya no incluye bajas automaticas"
# The automated audit extraction used to carry rows whose access was removed
# automatically for not meeting the 45-day login policy (I/J columns filled
# with "NO" / the policy message). This edit drops those auto-removal rows
# from PROCOTIZA and PRODESK so the sheets only keep the still-relevant
# records.

$wb = $excel.ActiveWorkbook

# --- PROCOTIZA: remove the 18 "baja automatica" rows (old rows 6-23),
#     keeping what used to be rows 24-25 (now rows 6-7). ---
$procotiza = $wb.Worksheets.Item("PROCOTIZA")
$procotiza.Rows("6:23").Delete() | Out-Null

# Column B (names) narrows slightly and column J (the old "motivo" free-text
# column) shrinks now that the long auto-removal message is no longer used.
$procotiza.Columns.Item(2).ColumnWidth = 35.166666666666664
$procotiza.Columns.Item(10).ColumnWidth = 14

# --- PRODESK: remove the single "baja automatica" row (old row 6, JAQUELINE
#     SORIA CHAVEZ), shifting the rest up. ---
$prodesk = $wb.Worksheets.Item("PRODESK")
$prodesk.Rows("6:6").Delete() | Out-Null

# Same column J narrowing as PROCOTIZA.
$prodesk.Columns.Item(10).ColumnWidth = 14
